# Apply cryptos list update (Fri Nov 29 04:55:37 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.385.46"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.577.25"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.406"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").Value = "3.573.59"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "4.259.90"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "96.674.57"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000258"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "3.571.87"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.492"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "511.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "3.769.27"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.96%  "
$ws.Range("E31").Value = "  +6.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "606.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("E40").Value = "  +6.60%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.905"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0417"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
